# EMPLOYEE_DATA.xlsx template update
# Commit: "column name changed in template. RM_ID AS RM_EMAIL_ID"
#
# The RM_NAME / RM_ID columns are renamed (RM_NAME -> RM_ID, RM_ID -> RM_EMAIL_ID)
# and the former "manager name" column (E) now holds the manager's e-mail address
# (with a mailto hyperlink, like column B already has for the employee's address).
# Two IS_ADMIN flags that were TRUE become FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row -------------------------------------------------------
$ws.Range("D1").Value = "RM_ID"
$ws.Range("E1").Value = "RM_EMAIL_ID"

# --- Column E: manager name -> manager e-mail (+ hyperlink) ----------
$rmEmails = @{
    2  = "SSC.NVEDAGIRI@CMA-CGM.COM"
    3  = "SSC.NVEDAGIRI@CMA-CGM.COM"
    4  = "SSC.NVEDAGIRI@CMA-CGM.COM"
    5  = "SSC.VRAMASAMY@CMA-CGM.COM"
    6  = "SSC.VRAMASAMY@CMA-CGM.COM"
    7  = "SSC.RM1@CMA-CGM.COM"
    8  = "SSC.RM1@CMA-CGM.COM"
    9  = "SSC.RM2@CMA-CGM.COM"
    10 = "SSC.RM2@CMA-CGM.COM"
    11 = "SSC.RM2@CMA-CGM.COM"
}

2..11 | ForEach-Object {
    $row = $_
    $cell = $ws.Range("E$row")
    $email = $rmEmails[$row]
    $cell.Value = $email
    $ws.Hyperlinks.Add($cell, "mailto:" + $email)
    $cell.Style = "Hyperlink"
}

# --- IS_ADMIN flags that flipped from TRUE to FALSE -------------------
$ws.Range("F2").Value = $false
$ws.Range("F9").Value = $false

# --- Selection moves from H13 to C13 -----------------------------------
$ws.Range("C13").Select()
